$wb = $excel.ActiveWorkbook

# 1. Remove 'bothered' from discard
$wsDiscard = $wb.Worksheets.Item("discard")
$wsDiscard.Rows.Item(2).Delete()

# 2. Add 'bothered' to negative (this creates the new "bot" shared string first)
$wsNeg = $wb.Worksheets.Item("negative")
$wsNeg.Cells.Item(27,1).Value = "bothered"
$wsNeg.Cells.Item(27,2).Value = "en"
$wsNeg.Cells.Item(27,3).Value = "bot"
$wsNeg.Range("C28").Select()

# 3. Update BAUM1 row in samples_retained (note text created after "bot")
$wsMain = $wb.Worksheets.Item("samples_retained")
$wsMain.Cells.Item(4,4).Value = 880
$wsMain.Cells.Item(4,8).Value = "contempt, surprise, unsure, and boredom mapped to negative; labels determined by interrater consensus; some of the mp4s might not have audio!; interest mapped to positive"

# 4. Final selection on samples_retained
$wsMain.Range("H5").Select()
